$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename header row: "_old" suffix -> "_FV2410", "_new" suffix -> "_FV2504"
#    (column 11 "diff" is left untouched)
$fv2410Headers = @("Segmentname_FV2410","Segmentgruppe_FV2410","Segment_FV2410","Datenelement_FV2410","Segment ID_FV2410","Code_FV2410","Qualifier_FV2410","Beschreibung_FV2410","Bedingungsausdruck_FV2410","Bedingung_FV2410")
$fv2504Headers = @("Segmentname_FV2504","Segmentgruppe_FV2504","Segment_FV2504","Datenelement_FV2504","Segment ID_FV2504","Code_FV2504","Qualifier_FV2504","Beschreibung_FV2504","Bedingungsausdruck_FV2504","Bedingung_FV2504")

for ($i = 0; $i -lt $fv2410Headers.Count; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $fv2410Headers[$i]
}

for ($i = 0; $i -lt $fv2504Headers.Count; $i++) {
    $ws.Cells.Item(1, $i + 12).Value = $fv2504Headers[$i]
}

# 2. Freeze the header row (top row frozen)
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
$ws.Range("A1").Select()

# 3. Convert the data range A1:U77 into an Excel Table named "Table1"
$rng = $ws.Range("A1:U77")
$tbl = $ws.ListObjects.Add(1, $rng, 0, 1)
$tbl.Name = "Table1"
